# Add knockout stage drawing: fill in the 8 team-name slots of the
# "Winner's Bracket" with the drawn team names, and clear the redundant
# "PrintYourBrackets.com" footer text (keeping its formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Team-name slots, one per bracket seed (1-8), typed in column B.
$teams = @(
    @{ Cell = "B4";  Name = "沙隆巴斯"; Family = 1 },
    @{ Cell = "B8";  Name = "金贝贝";   Family = 1 },
    @{ Cell = "B12"; Name = "约翰乔";   Family = 2 },
    @{ Cell = "B16"; Name = "阿土伯";   Family = 2 },
    @{ Cell = "B20"; Name = "莎拉公主"; Family = 1 },
    @{ Cell = "B24"; Name = "钱夫人";   Family = 1 },
    @{ Cell = "B28"; Name = "忍太郎";   Family = 1 },
    @{ Cell = "B32"; Name = "宫本宝藏"; Family = 1 }
)

foreach ($team in $teams) {
    $rng = $ws.Range($team.Cell)
    $rng.Value = $team.Name
    $rng.Font.Name = "微软雅黑"
    $rng.Font.Size = 8
    $rng.Font.Family = $team.Family
}

# The footer used to repeat the site's logo text under the empty branding
# cell above it (B50) -- remove the now-redundant duplicate, keeping the
# cell's existing border/alignment formatting intact.
$ws.Range("C63").Value = ""
